$d = $word.ActiveDocument

# 1. Title heading (appears twice: main H1 title and the bold "title" line near the end)
$d.Content.Find.Execute(
    "Play Colossal Gems Slot Free - Low Volatility and Giant Gem Symbols",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Play Colossal Gems Free: Review of Gameplay, Design, and Features", 2)

# 2. "What we like" bullet list items
$d.Content.Find.Execute(
    "Low volatility with small but frequent winnings",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Cartoon-like graphics appeal to slot players of all levels", 2)

# Item 3's old text is a substring of the italic meta-description sentence below it,
# so target only the specific bullet paragraph's Range to avoid touching that sentence.
$p50 = $d.Paragraphs.Item(51).Range
$p50.Find.Execute(
    "Respin and Free Spin features with giant gem symbols",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Low volatility provides frequent winnings", 2)

$d.Content.Find.Execute(
    "Autospin function with up to 500 autospins available",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Autospin function for convenient gameplay", 2)

$d.Content.Find.Execute(
    "Playable on virtually any device",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Compatible with all devices", 2)

# 3. "What we don't like" bullet list item
$d.Content.Find.Execute(
    "Graphics may not appeal to all players",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Progressive Jackpot Race is not available in all regions", 2)

# 4. Italic meta description sentence
$d.Content.Find.Execute(
    "Read our review of Colossal Gems slot game with low volatility Respin and Free Spin features with giant gem symbols. Play for free on any device.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Discover the colorful world of Colossal Gems and play this slot game for free. Learn about its gameplay, design, and bonus features.", 2)
